$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = 3
$ws.Range("F5").Value = -4
$ws.Range("F10").Value = -4
$ws.Range("F18").Value = 0
